$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Move the hidden "_GoBack" bookmark from the empty paragraph right
#    before "AREA DE TEORIA DE LA SENAL Y COMUNICACIONES" to wrap the
#    run containing "DEVELOPMENT OF RFID APPLICATIONS..." instead.
#    Bookmark names are unique, so adding a bookmark with the same
#    name moves it (the old one is implicitly removed).
# ------------------------------------------------------------------
$devRange = $null
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*DEVELOPMENT OF RFID APPLICATIONS*") {
        $devRange = $para.Range.Duplicate
        break
    }
}
if ($devRange -ne $null) {
    # Trim the trailing paragraph mark so the bookmark end stays inside
    # this paragraph (right after the text run), matching the target.
    $devRange.MoveEnd(1, -1) | Out-Null
    $d.Bookmarks.Add("_GoBack", $devRange) | Out-Null
}

# ------------------------------------------------------------------
# 2. Insert a new "COTUTOR: GUILLERMO ALVAREZ NARCIANDI" paragraph
#    right after the "TUTOR: YURI ALVAREZ LOPEZ" paragraph, matching
#    its formatting (it is duplicated from that paragraph mark).
# ------------------------------------------------------------------
$tutorPara = $null
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "TUTOR:*") {
        $tutorPara = $para
        break
    }
}
if ($tutorPara -ne $null) {
    $tutorPara.Range.InsertParagraphAfter() | Out-Null
    $newIndex = $tutorPara.Index + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $newPara.Range.Text = "COTUTOR: GUILLERMO ÁLVAREZ NARCIANDI"
}

# ------------------------------------------------------------------
# 3. Normalize a few built-in styles the same way Word itself does
#    when it resaves the package (adds uiPriority / unhideWhenUsed).
# ------------------------------------------------------------------
$normalize = @{
    "Fuentedeprrafopredeter" = 1
    "Tablanormal" = 99
    "Sinlista" = 99
}
foreach ($styleName in $normalize.Keys) {
    $style = $d.Styles.Item($styleName)
    if ($style -ne $null) {
        $style.Priority = $normalize[$styleName]
        $style.UnhideWhenUsed = $true
    }
}
